$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row correct-answer marks value: 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row correct marks value: 57 -> 95
$ws.Range("B12").Value = 95

# Update the correct/total marks summary text: 53/84 -> 95/140
$ws.Range("E12").Value = "95/140"
